$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 236 (shifts rows 236:315 down to 237:316,
# extending the used range to A1:R316).
$ws.Rows.Item(236).Insert()

# Populate the newly inserted row with its data.
$ws.Range("A236").Value = 3
$ws.Range("B236").Value = "Femacal de La Calera"
$ws.Range("C236").Value = "Coquimbo"
$ws.Range("D236").Value = 44524
$ws.Range("E236").Value = 5
$ws.Range("F236").Value = 100112021
$ws.Range("G236").Value = "Ají"
$ws.Range("H236").Value = "Inferno"
$ws.Range("I236").Value = "Primera"
$ws.Range("J236").Value = 73
$ws.Range("K236").Value = 22000
$ws.Range("L236").Value = 23000
$ws.Range("M236").Value = 22521
$ws.Range("N236").Value = "`$/caja 15 kilos"
$ws.Range("O236").Value = "Limache"
$ws.Range("P236").Value = 1501
$ws.Range("Q236").Value = 15
$ws.Range("R236").Value = "Hortaliza"
